$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.229.75'
$ws.Range("E2").Value = '  +11.35%  '

$ws.Range("D3").Value = '1.682.99'
$ws.Range("E3").Value = '  +7.04%  '

$ws.Range("D4").Value = '''1.004'
$ws.Range("E4").Value = '  +0.19%  '

$ws.Range("D5").Value = '''309.62'
$ws.Range("E5").Value = '  +8.34%  '

$ws.Range("D6").Value = '''0.9987'
$ws.Range("E6").Value = '  +1.44%  '

$ws.Range("D7").Value = '''0.3737'
$ws.Range("E7").Value = '  +1.45%  '

$ws.Range("D8").Value = '''0.3460'
$ws.Range("E8").Value = '  +5.59%  '

$ws.Range("D9").Value = '''47.87'
$ws.Range("E9").Value = '  +14.72%  '

$ws.Range("D10").Value = '''1.191'
$ws.Range("E10").Value = '  +5.05%  '

$ws.Range("D11").Value = '''0.07324'
$ws.Range("E11").Value = '  +4.07%  '

$ws.Range("D12").Value = '''1.000'
$ws.Range("E12").Value = '  +0.27%  '

$ws.Range("D13").Value = '''20.51'
$ws.Range("E13").Value = '  +3.00%  '

$ws.Range("D14").Value = '''6.135'
$ws.Range("E14").Value = '  +5.13%  '

$ws.Range("D15").Value = '''6.799'
$ws.Range("E15").Value = '  +4.37%  '

$ws.Range("D16").Value = '1.684.18'
$ws.Range("E16").Value = '  +7.35%  '

$ws.Range("D17").Value = '''0.00001112'
$ws.Range("E17").Value = '  +3.97%  '

$ws.Range("D18").Value = '''0.9984'
$ws.Range("E18").Value = '  +1.44%  '

$ws.Range("D19").Value = '''0.06727'
$ws.Range("E19").Value = '  +8.27%  '

$ws.Range("D20").Value = '''82.15'
$ws.Range("E20").Value = '  +10.24%  '

$ws.Range("E21").Value = '  +2.18%  '

$ws.Range("D22").Value = '''6.120'
$ws.Range("E22").Value = '  +4.41%  '

$ws.Range("D23").Value = '''12.06'
$ws.Range("E23").Value = '  +3.95%  '

$ws.Range("D24").Value = '24.209.48'
$ws.Range("E24").Value = '  +11.35%  '

$ws.Range("D25").Value = '''2.419'
$ws.Range("E25").Value = '  +2.06%  '

$ws.Range("D26").Value = '''2.683'
$ws.Range("E26").Value = '  +13.22%  '

$ws.Range("D27").Value = '''3.368'
$ws.Range("E27").Value = '  -9.02%  '

$ws.Range("D28").Value = '''153.93'
$ws.Range("E28").Value = '  +3.11%  '

$ws.Range("D29").Value = '''19.67'
$ws.Range("E29").Value = '  +7.80%  '

$ws.Range("D30").Value = '1.869.35'
$ws.Range("E30").Value = '  +7.34%  '

$ws.Range("D31").Value = '''127.31'
$ws.Range("E31").Value = '  +5.69%  '

$ws.Range("D32").Value = '''6.461'
$ws.Range("E32").Value = '  +18.91%  '

$ws.Range("D33").Value = '''4.108'
$ws.Range("E33").Value = '  +0.36%  '

$ws.Range("D34").Value = '''0.9886'
$ws.Range("E34").Value = '  +9.10%  '

$ws.Range("D35").Value = '''1.787'
$ws.Range("E35").Value = '  +11.95%  '

$ws.Range("D36").Value = '''0.08488'
$ws.Range("E36").Value = '  +3.56%  '

$ws.Range("D37").Value = '''12.48'
$ws.Range("E37").Value = '  +8.43%  '

$ws.Range("D38").Value = '''0.06490'
$ws.Range("E38").Value = '  +7.43%  '

$ws.Range("B39").Value = 'FraxShare'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D39").Value = '''8.969'
$ws.Range("E39").Value = '  +10.39%  '

$ws.Range("B40").Value = 'InternetComputer(DFINITY)'
$ws.Range("C40").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D40").Value = '''5.389'
$ws.Range("E40").Value = '  +5.35%  '

$ws.Range("D41").Value = '''0.02357'
$ws.Range("E41").Value = '  +8.72%  '

$ws.Range("D42").Value = '''1.278'
$ws.Range("E42").Value = '  +3.24%  '

$ws.Range("D43").Value = '''0.2136'
$ws.Range("E43").Value = '  +7.13%  '

$ws.Range("D44").Value = '''0.6218'
$ws.Range("E44").Value = '  +8.28%  '

$ws.Range("D45").Value = '''0.9977'
$ws.Range("E45").Value = '  +1.39%  '

$ws.Range("D46").Value = '''13.31'
$ws.Range("E46").Value = '  +3.75%  '

$ws.Range("D47").Value = '''3.810'
$ws.Range("E47").Value = '  +5.28%  '

$ws.Range("D48").Value = '''0.5983'
$ws.Range("E48").Value = '  +6.50%  '

$ws.Range("D49").Value = '''127.48'
$ws.Range("E49").Value = '  +2.28%  '

$ws.Range("D50").Value = '''2.043'
$ws.Range("E50").Value = '  +6.85%  '

$ws.Range("E51").Value = '  +6.86%  '
